$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.235.29'
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("E2").Style = $ws.Range("B2").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.830.18'
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E3").Style = $ws.Range("B3").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("E4").Style = $ws.Range("B4").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.14'
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("E5").Style = $ws.Range("B5").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6157'
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E6").Style = $ws.Range("B6").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E7").Style = $ws.Range("B7").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07347'
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("E8").Style = $ws.Range("B8").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2915'
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E9").Style = $ws.Range("B9").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.26'
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("E10").Style = $ws.Range("B10").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07659'
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("E11").Style = $ws.Range("B11").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.839.55'
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("E12").Style = $ws.Range("B12").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.987'
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("E13").Style = $ws.Range("B13").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6735'
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("E14").Style = $ws.Range("B14").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.57'
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("E15").Style = $ws.Range("B15").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008950'
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.15%  '
$ws.Range("E16").Style = $ws.Range("B16").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.874'
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("E17").Style = $ws.Range("B17").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.226.90'
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("E18").Style = $ws.Range("B18").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.095.30'
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("E19").Style = $ws.Range("B19").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.80'
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("E20").Style = $ws.Range("B20").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("E21").Style = $ws.Range("B21").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("E22").Style = $ws.Range("B22").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.388'
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.97%  '
$ws.Range("E23").Style = $ws.Range("B23").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.0000'
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E24").Style = $ws.Range("B24").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.74'
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("E25").Style = $ws.Range("B25").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.547'
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("E26").Style = $ws.Range("B26").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1391'
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("E27").Style = $ws.Range("B27").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.65'
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("E28").Style = $ws.Range("B28").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.494'
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E29").Style = $ws.Range("B29").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05775'
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.33%  '
$ws.Range("E30").Style = $ws.Range("B30").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.232'
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.00%  '
$ws.Range("E31").Style = $ws.Range("B31").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.088'
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("E32").Style = $ws.Range("B32").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.104'
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("E33").Style = $ws.Range("B33").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.858'
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E34").Style = $ws.Range("B34").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.136'
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E35").Style = $ws.Range("B35").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7213'
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.41%  '
$ws.Range("E36").Style = $ws.Range("B36").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.613'
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("E37").Style = $ws.Range("B37").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.862'
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.53%  '
$ws.Range("E38").Style = $ws.Range("B38").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.223.59'
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.09%  '
$ws.Range("E39").Style = $ws.Range("B39").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01764'
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("E40").Style = $ws.Range("B40").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.209'
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("E41").Style = $ws.Range("B41").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9067'
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("E42").Style = $ws.Range("B42").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("E43").Style = $ws.Range("B43").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.017.92'
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.21%  '
$ws.Range("E44").Style = $ws.Range("B44").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.87'
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E45").Style = $ws.Range("B45").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.58'
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("E46").Style = $ws.Range("B46").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5048'
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("E47").Style = $ws.Range("B47").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.12%  '
$ws.Range("E48").Style = $ws.Range("B48").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1182'
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.98%  '
$ws.Range("E49").Style = $ws.Range("B49").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.210'
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("E50").Style = $ws.Range("B50").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4041'
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.20%  '
$ws.Range("E51").Style = $ws.Range("B51").Style
